# Apply the commit's changes:
#  1. Revert the SiteName values in AntennaMetadata back to include
#     "Stationary Antenna" suffix (undo an earlier change that dropped it).
#  2. Restore the active cell selection on that sheet to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AntennaMetadata")
$ws.Activate()

# Rows 7-18 of column B hold the site names that need the suffix restored.
$rows = 7..18
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    $current = [string]$cell.Value2
    if ($current -notmatch "Stationary Antenna$") {
        $cell.Value = "$current Stationary Antenna"
    }
}

# Restore the previously-selected cell.
$ws.Range("B8").Select()
